$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Docente\(s\) Respons") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find 'Docente(s) Responsável(eis)' paragraph"
}

$headingPara = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new empty paragraph right after the heading.
$headingPara.Range.InsertParagraphAfter()

# That new empty paragraph is now the next one in the collection.
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Build the list-bullet paragraph containing the two professors, with a
# manual line break (not a new paragraph) between them, matching:
#   <w:r><w:t>5840730 - Antonio Jefferson da Silva Machado</w:t><w:br/></w:r>
#   <w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t></w:r>
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>5840730 - Antonio Jefferson da Silva Machado</w:t><w:br/></w:r>' +
       '<w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t></w:r>' +
       '</w:p>'

$newPara.Range.InsertXML($xml)
